$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New / updated data for rows 2-63 (columns A-D): A=index, B/C/D=measurements
$data = @(
    @(0, 0.8965421359345447, 7.298082327637384, 0.3793633423332122),
    @(1, 8.301818418868006, 10.44596042146457, 0.4499834365469383),
    @(2, 8.947782835694005, 7.91930493012805, 0.1247675139566222),
    @(3, 9.667468786393146, 4.796950672110532, 0.2700797743102612),
    @(4, 10.28511505101037, 8.152230305063135, 0.2287448513396756),
    @(5, 10.98187487424624, 3.046544535032166, 0.1086383738803665),
    @(6, 18.818786254039, 7.217635909333001, 0.5934006672536082),
    @(7, 23.52236998060614, 3.014971579036672, 0.2614872691718454),
    @(8, 25.48019593760126, 3.950647260058438, 0.08690586087957825),
    @(9, 26.5928910352749, 4.284385397189364, 0.1315994566268524),
    @(10, 27.49326239224975, 6.281582029500447, 0.2554409327148626),
    @(11, 28.67217840174301, 3.201656122060723, 0.4409916176385413),
    @(12, 29.0872638478475, 5.317349241289044, 0.3481971183012567),
    @(13, 29.13433449792829, 1.553110825352895, 0.3827122082263076),
    @(14, 33.31344962706875, 4.655294150950249, 0.2322501791772211),
    @(15, 33.46626526568765, 8.352433481669063, 0.3236750229241253),
    @(16, 34.02572465024981, 11.97941476008097, 0.3100179782393435),
    @(17, 36.3041776190417, 8.565226556325245, 0.3946103208931769),
    @(18, 37.65866264797165, 7.479571030753851, 0.1874095484715846),
    @(19, 38.10548982768861, 10.39563134747012, 0.3437712512348252),
    @(20, 39.77264523000445, 4.175622857689618, 0.3497724084131154),
    @(21, 43.20031993748078, 11.24009649528424, 0.4683312378745402),
    @(22, 44.54741048379763, 2.366973452469956, 0.2642572928778093),
    @(23, 51.26304501211248, 5.664425093474009, 0.1239531232626713),
    @(24, 56.17732117127385, 6.264133485491076, 0.5578644322875541),
    @(25, 58.24272937111937, 3.592640448893191, 0.06338052830646256),
    @(26, 60.32225634069512, 7.000308507356961, 0.2503395182458628),
    @(27, 61.40327556477762, 3.554957739801926, 0.3492449204960844),
    @(28, 63.71360383559794, 8.975236620318245, 0.1647345946937963),
    @(29, 63.77526595337387, 4.097891891907183, 0.1880639505245041),
    @(30, 64.40688592145388, 2.099561012953582, 0.3275549533842315),
    @(31, 65.21832139009115, 6.70856259791403, 0.1952624569355413),
    @(32, 66.71995400204393, 7.1703394078268, 0.2436535658656854),
    @(33, 68.91523201554011, 3.831608558016138, 0.1160557613118831),
    @(34, 69.13152289931138, 1.599041482900019, 0.1017988927188732),
    @(35, 69.983931224157, 2.36340673359018, 0.3998927538449091),
    @(36, 70.83518267337713, 7.036870762551326, 0.09968710125389274),
    @(37, 71.0925799158456, 5.502645871840715, 0.1456503631135305),
    @(38, 71.75356046258727, 6.178776118242983, 0.1827856557639106),
    @(39, 73.96119913114121, 7.002062737136058, 0.2163393333653746),
    @(40, 74.25667484109937, 5.660801761004527, 0.1711217976112475),
    @(41, 74.95352048484533, 6.033425315306188, 0.3743429270743347),
    @(42, 77.08736680772142, 2.42174107131811, 0.2867479474923859),
    @(43, 81.34664307608156, 7.218925660712713, 0.4940824646792574),
    @(44, 82.81693637586855, 5.066452928080667, 0.3976908213298642),
    @(45, 83.60154211229916, 4.725432490727738, 0.2745508906890788),
    @(46, 84.95037815857361, 6.321283319039207, 0.2372408024836736),
    @(47, 85.03032663603364, 4.617281666233332, 0.2729205658817148),
    @(48, 86.01175251103095, 10.88200036969712, 0.4118126931336475),
    @(49, 88.22207693682422, 2.887266122394368, 0.2550296144280473),
    @(50, 88.57435570029567, 6.512145844283001, 0.2515291654005461),
    @(51, 88.77953230338134, 5.463196103702572, 0.2724970337183749),
    @(52, 88.82572052857141, 1.195073832557399, 0.3865686048792937),
    @(53, 89.58705122668877, 7.235776631826199, 0.4113387269621717),
    @(54, 90.71992526648373, 4.696900735059289, 0.131162580236482),
    @(55, 91.20784857875849, 13.00601327863072, 0.2538352883080532),
    @(56, 91.74988232384278, 8.075549233350856, 0.3879230531443659),
    @(57, 95.22375525248776, 9.520392423245179, 0.3617930538542745),
    @(58, 95.84877420139662, 10.63566223379511, 0.3849580136699124),
    @(59, 96.35609221045229, 3.568308354398348, 0.2754575048962948),
    @(60, 97.38858239107668, 11.06933794918109, 0.388725925939909),
    @(61, 99.35901017263028, 8.102563093620835, 0.3341511500749703)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Newly added rows (52-63) need the same formatting (bold, centered, bordered)
# as the existing index column cells (A2:A51). Copy format from A2 onto A52:A63.
$ws.Range("A2").Copy()
$ws.Range("A52:A63").PasteSpecial(-4122)
$excel.CutCopyMode = $false
